$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 12
$ws.Range("I12").Value = 0.4959967553826868
$ws.Range("J12").Value = 0.01059582026591913
$ws.Range("K12").Value = -0.01514285199917103
$ws.Range("L12").Value = 1.663557790814084

# Row 13
$ws.Range("I13").Value = 0.4994216603121239
$ws.Range("J13").Value = 0.009152511499156205
$ws.Range("K13").Value = -0.1213644036507729
$ws.Range("L13").Value = 1.921334939992117

# Row 15
$ws.Range("I15").Value = 0.5169561186157083
$ws.Range("J15").Value = 0.01016064815776268
$ws.Range("K15").Value = -0.2810304976759443
$ws.Range("L15").Value = 1.845989007204349
